$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data set: username -> password pairs
$data = @(
    @("EM01", "sentrifugo"),
    @("EM02", "sentrifugo"),
    @("EM03", "sentrifugo"),
    @("EM04", "sentrifugo"),
    @("EM05", "sentrifugo"),
    @("EM06", "sentrifugo"),
    @("EM07", "sentrifugo"),
    @("EM08", "sentrifugo"),
    @("AGCY8", "sentrifugo"),
    @("US09", "sentrifugo")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Fit column B to its content
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Update selection to match the target state
$ws.Range("C10").Select() | Out-Null
